# Insert a new data row at row 61 (shifts existing rows 61-94 down to 62-95)
# and populate it with the new Papaya market-price record.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(61).Insert()

$ws.Range("A61").Value = 10
$ws.Range("B61").Value = "Vega Modelo de Temuco"
$ws.Range("C61").Value = "La Araucanía"
$ws.Range("D61").Value = 44981
$ws.Range("E61").Value = 9
$ws.Range("F61").Value = "Fruta"
$ws.Range("G61").Value = 100108
$ws.Range("H61").Value = "Tropicales y subtropicales"
$ws.Range("I61").Value = 100108004
$ws.Range("J61").Value = "Papaya"
$ws.Range("K61").Value = "Cultivar IV Región"
$ws.Range("L61").Value = "Primera"
$ws.Range("M61").Value = 20
$ws.Range("N61").Value = 40000
$ws.Range("O61").Value = 40000
$ws.Range("P61").Value = 40000
$ws.Range("Q61").Value = "$/caja 15 kilos granel"
$ws.Range("R61").Value = "Provincia de Limarí"
$ws.Range("S61").Value = 2667
$ws.Range("T61").Value = 15

Write-Host "Done. UsedRange:" $ws.UsedRange.Address()
